$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 2).Value = 0.1865889212827988
$ws.Cells.Item(2, 3).Value = 0.5685131195335277
$ws.Cells.Item(2, 10).Value = 0.02332361516034985
$ws.Cells.Item(2, 16).Value = 0.1239067055393586
$ws.Cells.Item(2, 19).Value = 0.09766763848396501

# Row 3
$ws.Cells.Item(3, 2).Value = 0.0072992700729927
$ws.Cells.Item(3, 3).Value = 0.0389294403892944
$ws.Cells.Item(3, 10).Value = 0.05109489051094891
$ws.Cells.Item(3, 16).Value = 0.7226277372262774
$ws.Cells.Item(3, 19).Value = 0.1800486618004866

# Row 4
$ws.Cells.Item(4, 10).Value = 0.04705882352941176
$ws.Cells.Item(4, 16).Value = 0.6588235294117647
$ws.Cells.Item(4, 19).Value = 0.2941176470588235

# Row 6
$ws.Cells.Item(6, 2).Value = 0.07692307692307693
$ws.Cells.Item(6, 4).Value = 0.0170940170940171
$ws.Cells.Item(6, 5).Value = 0.004273504273504274
$ws.Cells.Item(6, 6).Value = 0.08547008547008547
$ws.Cells.Item(6, 10).Value = 0.2542735042735043
$ws.Cells.Item(6, 15).Value = 0.0170940170940171
$ws.Cells.Item(6, 17).Value = 0.1773504273504274
$ws.Cells.Item(6, 18).Value = 0.07264957264957266
$ws.Cells.Item(6, 19).Value = 0.2948717948717949

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1313432835820895
$ws.Cells.Item(7, 4).Value = 0.01194029850746269
$ws.Cells.Item(7, 6).Value = 0.03880597014925373
$ws.Cells.Item(7, 10).Value = 0.1074626865671642
$ws.Cells.Item(7, 15).Value = 0.02686567164179104
$ws.Cells.Item(7, 17).Value = 0.217910447761194
$ws.Cells.Item(7, 18).Value = 0.09850746268656717
$ws.Cells.Item(7, 19).Value = 0.3671641791044776

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1140724946695096
$ws.Cells.Item(8, 4).Value = 0.0138592750533049
$ws.Cells.Item(8, 5).Value = 0.003198294243070362
$ws.Cells.Item(8, 6).Value = 0.07462686567164178
$ws.Cells.Item(8, 10).Value = 0.1108742004264392
$ws.Cells.Item(8, 15).Value = 0.01812366737739872
$ws.Cells.Item(8, 17).Value = 0.2025586353944563
$ws.Cells.Item(8, 18).Value = 0.09488272921108742
$ws.Cells.Item(8, 19).Value = 0.3678038379530917

# Row 9
$ws.Cells.Item(9, 2).Value = 0.1262886597938144
$ws.Cells.Item(9, 4).Value = 0.03092783505154639
$ws.Cells.Item(9, 6).Value = 0.09020618556701031
$ws.Cells.Item(9, 10).Value = 0.1082474226804124
$ws.Cells.Item(9, 15).Value = 0.01030927835051546
$ws.Cells.Item(9, 17).Value = 0.1701030927835052
$ws.Cells.Item(9, 18).Value = 0.1108247422680412
$ws.Cells.Item(9, 19).Value = 0.3530927835051547

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1160484759456482
$ws.Cells.Item(10, 4).Value = 0.01909658464928388
$ws.Cells.Item(10, 5).Value = 0.0007344840249724568
$ws.Cells.Item(10, 6).Value = 0.05728975394785164
$ws.Cells.Item(10, 10).Value = 0.1094381197208961
$ws.Cells.Item(10, 15).Value = 0.01505692251193537
$ws.Cells.Item(10, 17).Value = 0.2262210796915167
$ws.Cells.Item(10, 18).Value = 0.1053984575835476
$ws.Cells.Item(10, 19).Value = 0.3507161219243481

# Row 11
$ws.Cells.Item(11, 6).Value = 0.001886792452830189
$ws.Cells.Item(11, 7).Value = 0.1471698113207547
$ws.Cells.Item(11, 10).Value = 0.1018867924528302
$ws.Cells.Item(11, 11).Value = 0.2169811320754717
$ws.Cells.Item(11, 12).Value = 0.5150943396226415
$ws.Cells.Item(11, 19).Value = 0.0169811320754717

# Row 12
$ws.Cells.Item(12, 7).Value = 0.721830985915493
$ws.Cells.Item(12, 10).Value = 0.2042253521126761
$ws.Cells.Item(12, 11).Value = 0.0176056338028169
$ws.Cells.Item(12, 12).Value = 0.04577464788732395
$ws.Cells.Item(12, 19).Value = 0.01056338028169014

# Row 13
$ws.Cells.Item(13, 7).Value = 0.7045454545454546
$ws.Cells.Item(13, 10).Value = 0.2727272727272727
$ws.Cells.Item(13, 19).Value = 0.02272727272727273

# Row 15
$ws.Cells.Item(15, 6).Value = 0.02132196162046908
$ws.Cells.Item(15, 8).Value = 0.1513859275053305
$ws.Cells.Item(15, 9).Value = 0.05330490405117271
$ws.Cells.Item(15, 10).Value = 0.3795309168443497
$ws.Cells.Item(15, 11).Value = 0.06183368869936034
$ws.Cells.Item(15, 13).Value = 0.01066098081023454
$ws.Cells.Item(15, 14).Value = 0.002132196162046908
$ws.Cells.Item(15, 15).Value = 0.09594882729211088
$ws.Cells.Item(15, 19).Value = 0.2238805970149254

# Row 16
$ws.Cells.Item(16, 6).Value = 0.02358490566037736
$ws.Cells.Item(16, 8).Value = 0.1863207547169811
$ws.Cells.Item(16, 9).Value = 0.08018867924528301
$ws.Cells.Item(16, 10).Value = 0.3844339622641509
$ws.Cells.Item(16, 11).Value = 0.08490566037735849
$ws.Cells.Item(16, 13).Value = 0.02830188679245283
$ws.Cells.Item(16, 14).Value = 0.002358490566037736
$ws.Cells.Item(16, 15).Value = 0.09433962264150944
$ws.Cells.Item(16, 19).Value = 0.1155660377358491

# Row 17
$ws.Cells.Item(17, 6).Value = 0.0301556420233463
$ws.Cells.Item(17, 8).Value = 0.1692607003891051
$ws.Cells.Item(17, 9).Value = 0.08754863813229571
$ws.Cells.Item(17, 10).Value = 0.4523346303501946
$ws.Cells.Item(17, 11).Value = 0.07782101167315175
$ws.Cells.Item(17, 13).Value = 0.01750972762645914
$ws.Cells.Item(17, 15).Value = 0.0632295719844358
$ws.Cells.Item(17, 19).Value = 0.1021400778210117

# Row 18
$ws.Cells.Item(18, 6).Value = 0.03490759753593429
$ws.Cells.Item(18, 8).Value = 0.1765913757700205
$ws.Cells.Item(18, 9).Value = 0.09650924024640657
$ws.Cells.Item(18, 10).Value = 0.4353182751540041
$ws.Cells.Item(18, 11).Value = 0.08624229979466119
$ws.Cells.Item(18, 13).Value = 0.01848049281314168
$ws.Cells.Item(18, 15).Value = 0.06365503080082136
$ws.Cells.Item(18, 19).Value = 0.08829568788501027

# Row 19
$ws.Cells.Item(19, 6).Value = 0.01238134543953776
$ws.Cells.Item(19, 8).Value = 0.2199752373091209
$ws.Cells.Item(19, 9).Value = 0.07924061081304168
$ws.Cells.Item(19, 10).Value = 0.3962030540652084
$ws.Cells.Item(19, 11).Value = 0.09038382170862568
$ws.Cells.Item(19, 13).Value = 0.01981015270326042
$ws.Cells.Item(19, 14).Value = 0.001650846058605035
$ws.Cells.Item(19, 15).Value = 0.06397028477094512
$ws.Cells.Item(19, 19).Value = 0.116384647131655
